$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 2428.6  # H106: 2624.3125 -> 2428.6
$ws.Cells.Item(106, 9).Value = 2370.75  # I106: 2409.9092 -> 2370.75
$ws.Cells.Item(106, 10).Value = 2660  # J106: 3096 -> 2660
$ws.Cells.Item(106, 11).Value = 2370.75  # K106: 2409.9092 -> 2370.75
$ws.Cells.Item(106, 12).Value = 2660  # L106: 3096 -> 2660
$ws.Cells.Item(106, 13).Value = -1739.75  # M106: -1778.9092 -> -1739.75
$ws.Cells.Item(106, 14).Value = -3922  # N106: -4358 -> -3922

$ws.Cells.Item(132, 8).Value = 1967.2703  # H132: 2703.2917 -> 1967.2703
$ws.Cells.Item(132, 9).Value = 2068.4243  # I132: 2818.3635 -> 2068.4243
$ws.Cells.Item(132, 10).Value = 1132.75  # J132: 1437.5 -> 1132.75
$ws.Cells.Item(132, 11).Value = 6205.2729  # K132: 8455.0905 -> 6205.2729
$ws.Cells.Item(132, 12).Value = 3398.25  # L132: 4312.5 -> 3398.25
$ws.Cells.Item(132, 13).Value = -3675.2729  # M132: -5925.0905 -> -3675.2729
$ws.Cells.Item(132, 14).Value = -8458.25  # N132: -9372.5 -> -8458.25

$ws.Cells.Item(137, 8).Value = 11906201  # H137: 11906205 -> 11906201
$ws.Cells.Item(137, 9).Value = 1595.7142  # I137: 1547.9333 -> 1595.7142
$ws.Cells.Item(137, 10).Value = 35715412  # J137: 41667850 -> 35715412
$ws.Cells.Item(137, 11).Value = 4787.142599999999  # K137: 4643.7999 -> 4787.142599999999
$ws.Cells.Item(137, 12).Value = 107146236  # L137: 125003550 -> 107146236
$ws.Cells.Item(137, 13).Value = -2237.142599999999  # M137: -2093.7999 -> -2237.142599999999
$ws.Cells.Item(137, 14).Value = -107151336  # N137: -125008650 -> -107151336

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 75309.42999999999  # H2: 70308.266 -> 75309.42999999999
$ws.Cells.Item(2, 9).Value = 93393.91  # I2: 93511.27 -> 93393.91
$ws.Cells.Item(2, 10).Value = 8999.666999999999  # J2: 6500 -> 8999.666999999999
$ws.Cells.Item(2, 11).Value = 93393.91  # K2: 93511.27 -> 93393.91
$ws.Cells.Item(2, 12).Value = 8999.666999999999  # L2: 6500 -> 8999.666999999999
$ws.Cells.Item(2, 13).Value = -93280.91  # M2: -93398.27 -> -93280.91
$ws.Cells.Item(2, 14).Value = -9225.666999999999  # N2: -6726 -> -9225.666999999999

$ws.Cells.Item(32, 8).Value = 23892.941  # H32: 18588.984 -> 23892.941
$ws.Cells.Item(32, 9).Value = 25469.777  # I32: 18943.066 -> 25469.777
$ws.Cells.Item(32, 10).Value = 12066.667  # J32: 14340 -> 12066.667
$ws.Cells.Item(32, 11).Value = 25469.777  # K32: 18943.066 -> 25469.777
$ws.Cells.Item(32, 12).Value = 12066.667  # L32: 14340 -> 12066.667
$ws.Cells.Item(32, 13).Value = -25182.777  # M32: -18656.066 -> -25182.777
$ws.Cells.Item(32, 14).Value = -12640.667  # N32: -14914 -> -12640.667

$ws.Cells.Item(63, 8).Value = 3500  # H63: 2102.6086 -> 3500
$ws.Cells.Item(63, 9).Value = 5000  # I63: 2107.2727 -> 5000
$ws.Cells.Item(63, 11).Value = 5000  # K63: 2107.2727 -> 5000
$ws.Cells.Item(63, 13).Value = -4314  # M63: -1421.2727 -> -4314

$ws.Cells.Item(66, 8).Value = 3500  # H66: 2102.6086 -> 3500
$ws.Cells.Item(66, 9).Value = 5000  # I66: 2107.2727 -> 5000
$ws.Cells.Item(66, 11).Value = 25000  # K66: 10536.3635 -> 25000
$ws.Cells.Item(66, 13).Value = -21568  # M66: -7104.363499999999 -> -21568

$ws.Cells.Item(74, 8).Value = 813.7646999999999  # H74: 909.73914 -> 813.7646999999999
$ws.Cells.Item(74, 9).Value = 784  # I74: 1017.1111 -> 784
$ws.Cells.Item(74, 10).Value = 840.2222  # J74: 840.7143 -> 840.2222
$ws.Cells.Item(74, 11).Value = 784  # K74: 1017.1111 -> 784
$ws.Cells.Item(74, 12).Value = 840.2222  # L74: 840.7143 -> 840.2222
$ws.Cells.Item(74, 13).Value = 90  # M74: -143.1111 -> 90
$ws.Cells.Item(74, 14).Value = -2588.2222  # N74: -2588.7143 -> -2588.2222

$ws.Cells.Item(77, 8).Value = 813.7646999999999  # H77: 909.73914 -> 813.7646999999999
$ws.Cells.Item(77, 9).Value = 784  # I77: 1017.1111 -> 784
$ws.Cells.Item(77, 10).Value = 840.2222  # J77: 840.7143 -> 840.2222
$ws.Cells.Item(77, 11).Value = 3920  # K77: 5085.555499999999 -> 3920
$ws.Cells.Item(77, 12).Value = 4201.111  # L77: 4203.5715 -> 4201.111
$ws.Cells.Item(77, 13).Value = 448  # M77: -717.5554999999995 -> 448
$ws.Cells.Item(77, 14).Value = -12937.111  # N77: -12939.5715 -> -12937.111

$ws.Cells.Item(80, 8).Value = 15673.333  # H80: 15110 -> 15673.333
$ws.Cells.Item(80, 9).Value = 9800  # I80: 0 -> 9800
$ws.Cells.Item(80, 10).Value = 18610  # J80: 15110 -> 18610
$ws.Cells.Item(80, 11).Value = 9800  # K80: 0 -> 9800
$ws.Cells.Item(80, 12).Value = 18610  # L80: 15110 -> 18610
$ws.Cells.Item(80, 13).Value = -8802  # M80: None -> -8802
$ws.Cells.Item(80, 14).Value = -20606  # N80: -17106 -> -20606

$ws.Cells.Item(83, 8).Value = 15673.333  # H83: 15110 -> 15673.333
$ws.Cells.Item(83, 9).Value = 9800  # I83: 0 -> 9800
$ws.Cells.Item(83, 10).Value = 18610  # J83: 15110 -> 18610
$ws.Cells.Item(83, 11).Value = 29400  # K83: 0 -> 29400
$ws.Cells.Item(83, 12).Value = 55830  # L83: 45330 -> 55830
$ws.Cells.Item(83, 13).Value = -24408  # M83: None -> -24408
$ws.Cells.Item(83, 14).Value = -65814  # N83: -55314 -> -65814

$ws.Cells.Item(109, 8).Value = 0  # H109: 31333.334 -> 0
$ws.Cells.Item(109, 10).Value = 0  # J109: 31333.334 -> 0
$ws.Cells.Item(109, 12).Value = 0  # L109: 31333.334 -> 0
$ws.Cells.Item(109, 14).ClearContents()  # N109 was -34107.334

$ws.Cells.Item(112, 8).Value = 1673833.4  # H112: 1681600 -> 1673833.4
$ws.Cells.Item(112, 10).Value = 1673833.4  # J112: 1681600 -> 1673833.4
$ws.Cells.Item(112, 12).Value = 1673833.4  # L112: 1681600 -> 1673833.4
$ws.Cells.Item(112, 14).Value = -1676787.4  # N112: -1684554 -> -1676787.4

$ws.Cells.Item(116, 8).Value = 75309.42999999999  # H116: 70308.266 -> 75309.42999999999
$ws.Cells.Item(116, 9).Value = 93393.91  # I116: 93511.27 -> 93393.91
$ws.Cells.Item(116, 10).Value = 8999.666999999999  # J116: 6500 -> 8999.666999999999
$ws.Cells.Item(116, 11).Value = 93393.91  # K116: 93511.27 -> 93393.91
$ws.Cells.Item(116, 12).Value = 8999.666999999999  # L116: 6500 -> 8999.666999999999
$ws.Cells.Item(116, 13).Value = -91099.91  # M116: -91217.27 -> -91099.91
$ws.Cells.Item(116, 14).Value = -13587.667  # N116: -11088 -> -13587.667

$ws.Cells.Item(132, 8).Value = 3192.7083  # H132: 3653.4167 -> 3192.7083
$ws.Cells.Item(132, 9).Value = 2423  # I132: 2232.5 -> 2423
$ws.Cells.Item(132, 10).Value = 5501.8335  # J132: 4363.875 -> 5501.8335
$ws.Cells.Item(132, 11).Value = 7269  # K132: 6697.5 -> 7269
$ws.Cells.Item(132, 12).Value = 16505.5005  # L132: 13091.625 -> 16505.5005
$ws.Cells.Item(132, 13).Value = -4739  # M132: -4167.5 -> -4739
$ws.Cells.Item(132, 14).Value = -21565.5005  # N132: -18151.625 -> -21565.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 75309.42999999999  # H3: 70308.266 -> 75309.42999999999
$ws.Cells.Item(3, 9).Value = 93393.91  # I3: 93511.27 -> 93393.91
$ws.Cells.Item(3, 10).Value = 8999.666999999999  # J3: 6500 -> 8999.666999999999
$ws.Cells.Item(3, 11).Value = 93393.91  # K3: 93511.27 -> 93393.91
$ws.Cells.Item(3, 12).Value = 8999.666999999999  # L3: 6500 -> 8999.666999999999
$ws.Cells.Item(3, 13).Value = -93279.91  # M3: -93397.27 -> -93279.91
$ws.Cells.Item(3, 14).Value = -9227.666999999999  # N3: -6728 -> -9227.666999999999

$ws.Cells.Item(82, 8).Value = 15695.223  # H82: 16530.111 -> 15695.223
$ws.Cells.Item(82, 9).Value = 10451.4  # I82: 12154.2 -> 10451.4
$ws.Cells.Item(82, 10).Value = 22250  # J82: 22000 -> 22250
$ws.Cells.Item(82, 11).Value = 10451.4  # K82: 12154.2 -> 10451.4
$ws.Cells.Item(82, 12).Value = 22250  # L82: 22000 -> 22250
$ws.Cells.Item(82, 13).Value = -10068.4  # M82: -11771.2 -> -10068.4
$ws.Cells.Item(82, 14).Value = -23016  # N82: -22766 -> -23016

$ws.Cells.Item(85, 8).Value = 15695.223  # H85: 16530.111 -> 15695.223
$ws.Cells.Item(85, 9).Value = 10451.4  # I85: 12154.2 -> 10451.4
$ws.Cells.Item(85, 10).Value = 22250  # J85: 22000 -> 22250
$ws.Cells.Item(85, 11).Value = 10451.4  # K85: 12154.2 -> 10451.4
$ws.Cells.Item(85, 12).Value = 22250  # L85: 22000 -> 22250
$ws.Cells.Item(85, 13).Value = -9125.4  # M85: -10828.2 -> -9125.4
$ws.Cells.Item(85, 14).Value = -24902  # N85: -24652 -> -24902

$ws.Cells.Item(119, 8).Value = 17000  # H119: 19500 -> 17000
$ws.Cells.Item(119, 10).Value = 17000  # J119: 19500 -> 17000
$ws.Cells.Item(119, 12).Value = 17000  # L119: 19500 -> 17000
$ws.Cells.Item(119, 14).Value = -26676  # N119: -29176 -> -26676

$ws.Cells.Item(138, 8).Value = 43666.668  # H138: 0 -> 43666.668
$ws.Cells.Item(138, 10).Value = 43666.668  # J138: 0 -> 43666.668
$ws.Cells.Item(138, 12).Value = 43666.668  # L138: 0 -> 43666.668
$ws.Cells.Item(138, 14).Value = -53946.668  # N138: None -> -53946.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2166.8235  # H31: 2007.8948 -> 2166.8235
$ws.Cells.Item(31, 9).Value = 1683.2632  # I31: 1544.8182 -> 1683.2632
$ws.Cells.Item(31, 10).Value = 2779.3333  # J31: 2644.625 -> 2779.3333
$ws.Cells.Item(31, 11).Value = 1683.2632  # K31: 1544.8182 -> 1683.2632
$ws.Cells.Item(31, 12).Value = 2779.3333  # L31: 2644.625 -> 2779.3333
$ws.Cells.Item(31, 13).Value = -1388.2632  # M31: -1249.8182 -> -1388.2632
$ws.Cells.Item(31, 14).Value = -3369.3333  # N31: -3234.625 -> -3369.3333

$ws.Cells.Item(34, 8).Value = 2166.8235  # H34: 2007.8948 -> 2166.8235
$ws.Cells.Item(34, 9).Value = 1683.2632  # I34: 1544.8182 -> 1683.2632
$ws.Cells.Item(34, 10).Value = 2779.3333  # J34: 2644.625 -> 2779.3333
$ws.Cells.Item(34, 11).Value = 1683.2632  # K34: 1544.8182 -> 1683.2632
$ws.Cells.Item(34, 12).Value = 2779.3333  # L34: 2644.625 -> 2779.3333
$ws.Cells.Item(34, 13).Value = -1481.2632  # M34: -1342.8182 -> -1481.2632
$ws.Cells.Item(34, 14).Value = -3183.3333  # N34: -3048.625 -> -3183.3333

$ws.Cells.Item(68, 8).Value = 20780  # H68: 20977.777 -> 20780
$ws.Cells.Item(68, 10).Value = 20780  # J68: 20977.777 -> 20780
$ws.Cells.Item(68, 12).Value = 20780  # L68: 20977.777 -> 20780
$ws.Cells.Item(68, 14).Value = -22278  # N68: -22475.777 -> -22278

$ws.Cells.Item(71, 8).Value = 20780  # H71: 20977.777 -> 20780
$ws.Cells.Item(71, 10).Value = 20780  # J71: 20977.777 -> 20780
$ws.Cells.Item(71, 12).Value = 62340  # L71: 62933.33099999999 -> 62340
$ws.Cells.Item(71, 14).Value = -69828  # N71: -70421.33099999999 -> -69828

$ws.Cells.Item(74, 8).Value = 30000  # H74: 0 -> 30000
$ws.Cells.Item(74, 10).Value = 30000  # J74: 0 -> 30000
$ws.Cells.Item(74, 12).Value = 30000  # L74: 0 -> 30000
$ws.Cells.Item(74, 14).Value = -31748  # N74: None -> -31748

$ws.Cells.Item(77, 8).Value = 30000  # H77: 0 -> 30000
$ws.Cells.Item(77, 10).Value = 30000  # J77: 0 -> 30000
$ws.Cells.Item(77, 12).Value = 90000  # L77: 0 -> 90000
$ws.Cells.Item(77, 14).Value = -98736  # N77: None -> -98736

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 233.1  # H2: 204.25 -> 233.1
$ws.Cells.Item(2, 9).Value = 186.14285  # I2: 158.11111 -> 186.14285
$ws.Cells.Item(2, 11).Value = 1116.8571  # K2: 948.66666 -> 1116.8571
$ws.Cells.Item(2, 13).Value = -1003.8571  # M2: -835.66666 -> -1003.8571

$ws.Cells.Item(3, 8).Value = 6324.5835  # H3: 6348.8 -> 6324.5835
$ws.Cells.Item(3, 9).Value = 1975  # I3: 2525.5557 -> 1975
$ws.Cells.Item(3, 11).Value = 5925  # K3: 7576.6671 -> 5925
$ws.Cells.Item(3, 13).Value = -5813  # M3: -7464.6671 -> -5813

$ws.Cells.Item(107, 8).Value = 914.70734  # H107: 904.61365 -> 914.70734
$ws.Cells.Item(107, 9).Value = 275.85715  # I107: 289.3846 -> 275.85715
$ws.Cells.Item(107, 10).Value = 1245.963  # J107: 1162.6129 -> 1245.963
$ws.Cells.Item(107, 11).Value = 827.5714499999999  # K107: 868.1537999999999 -> 827.5714499999999
$ws.Cells.Item(107, 12).Value = 3737.889  # L107: 3487.8387 -> 3737.889
$ws.Cells.Item(107, 13).Value = 1092.42855  # M107: 1051.8462 -> 1092.42855
$ws.Cells.Item(107, 14).Value = -7577.889  # N107: -7327.8387 -> -7577.889

$ws.Cells.Item(113, 8).Value = 390.83  # H113: 400.71 -> 390.83
$ws.Cells.Item(113, 9).Value = 362.5  # I113: 356.4 -> 362.5
$ws.Cells.Item(113, 10).Value = 399.7763  # J113: 411.7875 -> 399.7763
$ws.Cells.Item(113, 11).Value = 1087.5  # K113: 1069.2 -> 1087.5
$ws.Cells.Item(113, 12).Value = 1199.3289  # L113: 1235.3625 -> 1199.3289
$ws.Cells.Item(113, 13).Value = 1082.5  # M113: 1100.8 -> 1082.5
$ws.Cells.Item(113, 14).Value = -5539.3289  # N113: -5575.3625 -> -5539.3289

$ws.Cells.Item(131, 8).Value = 15306.261  # H131: 15931.643 -> 15306.261
$ws.Cells.Item(131, 9).Value = 67804.664  # I131: 72620 -> 67804.664
$ws.Cells.Item(131, 10).Value = 1729.0862  # J131: 1759.5536 -> 1729.0862
$ws.Cells.Item(131, 11).Value = 203413.992  # K131: 217860 -> 203413.992
$ws.Cells.Item(131, 12).Value = 5187.2586  # L131: 5278.6608 -> 5187.2586
$ws.Cells.Item(131, 13).Value = -198373.992  # M131: -212820 -> -198373.992
$ws.Cells.Item(131, 14).Value = -15267.2586  # N131: -15358.6608 -> -15267.2586

$ws.Cells.Item(134, 8).Value = 5327.2104  # H134: 6383.5264 -> 5327.2104
$ws.Cells.Item(134, 9).Value = 4093.6155  # I134: 4928.7 -> 4093.6155
$ws.Cells.Item(134, 11).Value = 12280.8465  # K134: 14786.1 -> 12280.8465
$ws.Cells.Item(134, 13).Value = -7210.8465  # M134: -9716.099999999999 -> -7210.8465

$ws.Cells.Item(139, 8).Value = 1196.2667  # H139: 1253.3125 -> 1196.2667
$ws.Cells.Item(139, 9).Value = 541  # I139: 528.1818 -> 541
$ws.Cells.Item(139, 10).Value = 2998.25  # J139: 2848.6 -> 2998.25
$ws.Cells.Item(139, 11).Value = 1623  # K139: 1584.5454 -> 1623
$ws.Cells.Item(139, 12).Value = 8994.75  # L139: 8545.799999999999 -> 8994.75
$ws.Cells.Item(139, 13).Value = 3517  # M139: 3555.4546 -> 3517
$ws.Cells.Item(139, 14).Value = -19274.75  # N139: -18825.8 -> -19274.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4342.2383  # H70: 4139.4736 -> 4342.2383
$ws.Cells.Item(70, 9).Value = 4090.6428  # I70: 3746 -> 4090.6428
$ws.Cells.Item(70, 10).Value = 4845.4287  # J70: 4493.6 -> 4845.4287
$ws.Cells.Item(70, 11).Value = 4090.6428  # K70: 3746 -> 4090.6428
$ws.Cells.Item(70, 12).Value = 4845.4287  # L70: 4493.6 -> 4845.4287
$ws.Cells.Item(70, 13).Value = -3820.6428  # M70: -3476 -> -3820.6428
$ws.Cells.Item(70, 14).Value = -5385.4287  # N70: -5033.6 -> -5385.4287

$ws.Cells.Item(73, 8).Value = 4342.2383  # H73: 4139.4736 -> 4342.2383
$ws.Cells.Item(73, 9).Value = 4090.6428  # I73: 3746 -> 4090.6428
$ws.Cells.Item(73, 10).Value = 4845.4287  # J73: 4493.6 -> 4845.4287
$ws.Cells.Item(73, 11).Value = 4090.6428  # K73: 3746 -> 4090.6428
$ws.Cells.Item(73, 12).Value = 4845.4287  # L73: 4493.6 -> 4845.4287
$ws.Cells.Item(73, 13).Value = -3154.6428  # M73: -2810 -> -3154.6428
$ws.Cells.Item(73, 14).Value = -6717.4287  # N73: -6365.6 -> -6717.4287

$ws.Cells.Item(80, 8).Value = 3793.5  # H80: 3673.9333 -> 3793.5
$ws.Cells.Item(80, 10).Value = 3000  # J80: 2500 -> 3000
$ws.Cells.Item(80, 12).Value = 3000  # L80: 2500 -> 3000
$ws.Cells.Item(80, 14).Value = -4996  # N80: -4496 -> -4996

$ws.Cells.Item(83, 8).Value = 3793.5  # H83: 3673.9333 -> 3793.5
$ws.Cells.Item(83, 10).Value = 3000  # J83: 2500 -> 3000
$ws.Cells.Item(83, 12).Value = 15000  # L83: 12500 -> 15000
$ws.Cells.Item(83, 14).Value = -24984  # N83: -22484 -> -24984

$ws.Cells.Item(132, 8).Value = 2875.0625  # H132: 2928.375 -> 2875.0625
$ws.Cells.Item(132, 9).Value = 2284.3333  # I132: 2423.0435 -> 2284.3333
$ws.Cells.Item(132, 10).Value = 4647.25  # J132: 4219.778 -> 4647.25
$ws.Cells.Item(132, 11).Value = 6852.999899999999  # K132: 7269.130500000001 -> 6852.999899999999
$ws.Cells.Item(132, 12).Value = 13941.75  # L132: 12659.334 -> 13941.75
$ws.Cells.Item(132, 13).Value = -4322.999899999999  # M132: -4739.130500000001 -> -4322.999899999999
$ws.Cells.Item(132, 14).Value = -19001.75  # N132: -17719.334 -> -19001.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(110, 8).Value = 644  # H110: 29800 -> 644
$ws.Cells.Item(110, 10).Value = 644  # J110: 29800 -> 644
$ws.Cells.Item(110, 12).Value = 644  # L110: 29800 -> 644
$ws.Cells.Item(110, 14).Value = -8824  # N110: -37980 -> -8824

$ws.Cells.Item(132, 8).Value = 2780390.8  # H132: 2780487.2 -> 2780390.8
$ws.Cells.Item(132, 9).Value = 3970453.5  # I132: 4632176.5 -> 3970453.5
$ws.Cells.Item(132, 10).Value = 3578.111  # J132: 2953.1667 -> 3578.111
$ws.Cells.Item(132, 11).Value = 11911360.5  # K132: 13896529.5 -> 11911360.5
$ws.Cells.Item(132, 12).Value = 10734.333  # L132: 8859.500100000001 -> 10734.333
$ws.Cells.Item(132, 13).Value = -11908830.5  # M132: -13893999.5 -> -11908830.5
$ws.Cells.Item(132, 14).Value = -15794.333  # N132: -13919.5001 -> -15794.333

$ws.Cells.Item(136, 8).Value = 3478.5  # H136: 3453.5 -> 3478.5
$ws.Cells.Item(136, 9).Value = 1124.0769  # I136: 1083.0714 -> 1124.0769
$ws.Cells.Item(136, 10).Value = 9600  # J136: 11750 -> 9600
$ws.Cells.Item(136, 11).Value = 3372.2307  # K136: 3249.2142 -> 3372.2307
$ws.Cells.Item(136, 12).Value = 28800  # L136: 35250 -> 28800
$ws.Cells.Item(136, 13).Value = -822.2307000000001  # M136: -699.2142000000003 -> -822.2307000000001
$ws.Cells.Item(136, 14).Value = -33900  # N136: -40350 -> -33900

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 67164.75  # H46: 57097.832 -> 67164.75
$ws.Cells.Item(46, 10).Value = 67164.75  # J46: 57097.832 -> 67164.75
$ws.Cells.Item(46, 12).Value = 67164.75  # L46: 57097.832 -> 67164.75
$ws.Cells.Item(46, 14).Value = -67626.75  # N46: -57559.832 -> -67626.75

$ws.Cells.Item(62, 8).Value = 3439.8572  # H62: 3650 -> 3439.8572
$ws.Cells.Item(62, 9).Value = 3369.75  # I62: 3700 -> 3369.75
$ws.Cells.Item(62, 10).Value = 3533.3333  # J62: 3600 -> 3533.3333
$ws.Cells.Item(62, 11).Value = 3369.75  # K62: 3700 -> 3369.75
$ws.Cells.Item(62, 12).Value = 3533.3333  # L62: 3600 -> 3533.3333
$ws.Cells.Item(62, 13).Value = -2745.75  # M62: -3076 -> -2745.75
$ws.Cells.Item(62, 14).Value = -4781.3333  # N62: -4848 -> -4781.3333

$ws.Cells.Item(65, 8).Value = 3439.8572  # H65: 3650 -> 3439.8572
$ws.Cells.Item(65, 9).Value = 3369.75  # I65: 3700 -> 3369.75
$ws.Cells.Item(65, 10).Value = 3533.3333  # J65: 3600 -> 3533.3333
$ws.Cells.Item(65, 11).Value = 16848.75  # K65: 18500 -> 16848.75
$ws.Cells.Item(65, 12).Value = 17666.6665  # L65: 18000 -> 17666.6665
$ws.Cells.Item(65, 13).Value = -13728.75  # M65: -15380 -> -13728.75
$ws.Cells.Item(65, 14).Value = -23906.6665  # N65: -24240 -> -23906.6665

$ws.Cells.Item(81, 8).Value = 2250  # H81: 1281.8334 -> 2250
$ws.Cells.Item(81, 9).Value = 1500  # I81: 1284.75 -> 1500
$ws.Cells.Item(81, 10).Value = 3000  # J81: 1276 -> 3000
$ws.Cells.Item(81, 11).Value = 3000  # K81: 2569.5 -> 3000
$ws.Cells.Item(81, 12).Value = 6000  # L81: 2552 -> 6000
$ws.Cells.Item(81, 13).Value = -1939  # M81: -1508.5 -> -1939
$ws.Cells.Item(81, 14).Value = -8122  # N81: -4674 -> -8122

$ws.Cells.Item(84, 8).Value = 2250  # H84: 1281.8334 -> 2250
$ws.Cells.Item(84, 9).Value = 1500  # I84: 1284.75 -> 1500
$ws.Cells.Item(84, 10).Value = 3000  # J84: 1276 -> 3000
$ws.Cells.Item(84, 11).Value = 15000  # K84: 12847.5 -> 15000
$ws.Cells.Item(84, 12).Value = 30000  # L84: 12760 -> 30000
$ws.Cells.Item(84, 13).Value = -9696  # M84: -7543.5 -> -9696
$ws.Cells.Item(84, 14).Value = -40608  # N84: -23368 -> -40608

$ws.Cells.Item(134, 8).Value = 67164.75  # H134: 57097.832 -> 67164.75
$ws.Cells.Item(134, 10).Value = 67164.75  # J134: 57097.832 -> 67164.75
$ws.Cells.Item(134, 12).Value = 201494.25  # L134: 171293.496 -> 201494.25
$ws.Cells.Item(134, 14).Value = -206564.25  # N134: -176363.496 -> -206564.25

$ws.Cells.Item(136, 8).Value = 2356.2856  # H136: 2263.0454 -> 2356.2856
$ws.Cells.Item(136, 9).Value = 2501.625  # I136: 2338.423 -> 2501.625
$ws.Cells.Item(136, 10).Value = 2162.5  # J136: 2154.1667 -> 2162.5
$ws.Cells.Item(136, 11).Value = 7504.875  # K136: 7015.268999999999 -> 7504.875
$ws.Cells.Item(136, 12).Value = 6487.5  # L136: 6462.500100000001 -> 6487.5
$ws.Cells.Item(136, 13).Value = -4954.875  # M136: -4465.268999999999 -> -4954.875
$ws.Cells.Item(136, 14).Value = -11587.5  # N136: -11562.5001 -> -11587.5
